# Auto-generated edit script: updates Yojimbo_Profits market-data values
# across ALC/ARM/CRP/CUL/GSM/LTW/WVR sheets per the scheduled-runner sync.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 13
$ws.Range("H13").Value = 5050
$ws.Range("J13").Value = 5050
$ws.Range("L13").Value = 5050
$ws.Range("N13").Value = -5388

# Row 111
$ws.Range("H111").Value = 268.5
$ws.Range("I111").Value = 215.8
$ws.Range("J111").Value = 532
$ws.Range("K111").Value = 647.4000000000001
$ws.Range("L111").Value = 1596
$ws.Range("M111").Value = 2419.6
$ws.Range("N111").Value = -7730

# Row 132
$ws.Range("H132").Value = 1454.5814
$ws.Range("I132").Value = 1464.4147
$ws.Range("J132").Value = 1253
$ws.Range("K132").Value = 4393.2441
$ws.Range("L132").Value = 3759
$ws.Range("M132").Value = -1863.2441
$ws.Range("N132").Value = -8819

# Row 141
$ws.Range("H141").Value = 2217.756
$ws.Range("I141").Value = 1893.9429
$ws.Range("J141").Value = 4106.6665
$ws.Range("K141").Value = 5681.8287
$ws.Range("L141").Value = 12319.9995
$ws.Range("M141").Value = -501.8287
$ws.Range("N141").Value = -22679.9995

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 1047.9231
$ws.Range("I61").Value = 797.90247
$ws.Range("J61").Value = 1979.8182
$ws.Range("K61").Value = 797.90247
$ws.Range("L61").Value = 1979.8182
$ws.Range("M61").Value = -585.90247
$ws.Range("N61").Value = -2403.8182

# Row 122
$ws.Range("H122").Value = 8336022
$ws.Range("I122").Value = 16669409
$ws.Range("J122").Value = 2635.6667
$ws.Range("K122").Value = 50008227
$ws.Range("L122").Value = 7907.000100000001
$ws.Range("M122").Value = -50005777
$ws.Range("N122").Value = -12807.0001

# Row 136
$ws.Range("H136").Value = 1047.9231
$ws.Range("I136").Value = 797.90247
$ws.Range("J136").Value = 1979.8182
$ws.Range("K136").Value = 2393.70741
$ws.Range("L136").Value = 5939.4546
$ws.Range("M136").Value = 156.29259
$ws.Range("N136").Value = -11039.4546

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 833.3333
$ws.Range("I16").Value = 900
$ws.Range("J16").Value = 800
$ws.Range("K16").Value = 900
$ws.Range("L16").Value = 800
$ws.Range("M16").Value = -613
$ws.Range("N16").Value = -1374

# Row 31
$ws.Range("H31").Value = 32050.953
$ws.Range("I31").Value = 41386.152
$ws.Range("J31").Value = 16881.25
$ws.Range("K31").Value = 41386.152
$ws.Range("L31").Value = 16881.25
$ws.Range("M31").Value = -41091.152
$ws.Range("N31").Value = -17471.25

# Row 34
$ws.Range("H34").Value = 32050.953
$ws.Range("I34").Value = 41386.152
$ws.Range("J34").Value = 16881.25
$ws.Range("K34").Value = 41386.152
$ws.Range("L34").Value = 16881.25
$ws.Range("M34").Value = -41184.152
$ws.Range("N34").Value = -17285.25

# Row 58
$ws.Range("H58").Value = 1105.9524
$ws.Range("I58").Value = 1021.2407
$ws.Range("J58").Value = 1614.2222
$ws.Range("K58").Value = 1021.2407
$ws.Range("L58").Value = 1614.2222
$ws.Range("M58").Value = -818.2406999999999
$ws.Range("N58").Value = -2020.2222

# Row 62
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 2300
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 2300
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -3548

# Row 65
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 2300
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 11500
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -17740

# Row 113
$ws.Range("H113").Value = 833.3333
$ws.Range("I113").Value = 900
$ws.Range("J113").Value = 800
$ws.Range("K113").Value = 900
$ws.Range("L113").Value = 800
$ws.Range("M113").Value = 1270
$ws.Range("N113").Value = -5140

# Row 134
$ws.Range("H134").Value = 3110.8
$ws.Range("I134").Value = 2617.7144
$ws.Range("J134").Value = 10014
$ws.Range("K134").Value = 7853.1432
$ws.Range("L134").Value = 30042
$ws.Range("M134").Value = -5318.1432
$ws.Range("N134").Value = -35112

# Row 136
$ws.Range("H136").Value = 1105.9524
$ws.Range("I136").Value = 1021.2407
$ws.Range("J136").Value = 1614.2222
$ws.Range("K136").Value = 3063.7221
$ws.Range("L136").Value = 4842.6666
$ws.Range("M136").Value = -513.7221
$ws.Range("N136").Value = -9942.6666

$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 915.48486
$ws.Range("I113").Value = 1267.3889
$ws.Range("J113").Value = 493.2
$ws.Range("K113").Value = 3802.1667
$ws.Range("L113").Value = 1479.6
$ws.Range("M113").Value = -1632.1667
$ws.Range("N113").Value = -5819.6

# Row 114
$ws.Range("H114").Value = 393.07144
$ws.Range("I114").Value = 286.33334
$ws.Range("J114").Value = 473.125
$ws.Range("K114").Value = 859.0000200000001
$ws.Range("L114").Value = 1419.375
$ws.Range("M114").Value = 2394.99998
$ws.Range("N114").Value = -7927.375

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 2429.8572
$ws.Range("I80").Value = 2552.3333
$ws.Range("J80").Value = 2338
$ws.Range("K80").Value = 2552.3333
$ws.Range("L80").Value = 2338
$ws.Range("M80").Value = -1554.3333
$ws.Range("N80").Value = -4334

# Row 83
$ws.Range("H83").Value = 2429.8572
$ws.Range("I83").Value = 2552.3333
$ws.Range("J83").Value = 2338
$ws.Range("K83").Value = 12761.6665
$ws.Range("L83").Value = 11690
$ws.Range("M83").Value = -7769.666499999999
$ws.Range("N83").Value = -21674

# Row 97
$ws.Range("H97").Value = 1229.0938
$ws.Range("I97").Value = 1171.0435
$ws.Range("J97").Value = 1377.4445
$ws.Range("K97").Value = 1171.0435
$ws.Range("L97").Value = 1377.4445
$ws.Range("M97").Value = -675.0435
$ws.Range("N97").Value = -2369.4445

# Row 132
$ws.Range("H132").Value = 1243.3334
$ws.Range("I132").Value = 1276.279
$ws.Range("K132").Value = 3828.837
$ws.Range("M132").Value = -1298.837

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 1876.8667
$ws.Range("I7").Value = 1216.6
$ws.Range("J7").Value = 2207
$ws.Range("K7").Value = 1216.6
$ws.Range("L7").Value = 2207
$ws.Range("M7").Value = -1104.6
$ws.Range("N7").Value = -2431

# Row 82
$ws.Range("H82").Value = 2396.8928
$ws.Range("I82").Value = 1026.5
$ws.Range("J82").Value = 2945.05
$ws.Range("K82").Value = 1026.5
$ws.Range("L82").Value = 2945.05
$ws.Range("M82").Value = -665.5
$ws.Range("N82").Value = -3667.05

# Row 85
$ws.Range("H85").Value = 2396.8928
$ws.Range("I85").Value = 1026.5
$ws.Range("J85").Value = 2945.05
$ws.Range("K85").Value = 1026.5
$ws.Range("L85").Value = 2945.05
$ws.Range("M85").Value = 221.5
$ws.Range("N85").Value = -5441.05

# Row 126
$ws.Range("H126").Value = 1876.8667
$ws.Range("I126").Value = 1216.6
$ws.Range("J126").Value = 2207
$ws.Range("K126").Value = 3649.8
$ws.Range("L126").Value = 6621
$ws.Range("M126").Value = -1179.8
$ws.Range("N126").Value = -11561

# Row 132
$ws.Range("H132").Value = 1906.2373
$ws.Range("I132").Value = 1669.8541
$ws.Range("J132").Value = 2937.7273
$ws.Range("K132").Value = 5009.5623
$ws.Range("L132").Value = 8813.1819
$ws.Range("M132").Value = -2479.5623
$ws.Range("N132").Value = -13873.1819

# Row 136
$ws.Range("H136").Value = 1971.5316
$ws.Range("I136").Value = 1377.5264
$ws.Range("J136").Value = 3510.5454
$ws.Range("K136").Value = 4132.5792
$ws.Range("L136").Value = 10531.6362
$ws.Range("M136").Value = -1582.5792
$ws.Range("N136").Value = -15631.6362

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 1331.909
$ws.Range("I81").Value = 1294.5555
$ws.Range("J81").Value = 1500
$ws.Range("K81").Value = 2589.111
$ws.Range("L81").Value = 3000
$ws.Range("M81").Value = -1528.111
$ws.Range("N81").Value = -5122

# Row 84
$ws.Range("H84").Value = 1331.909
$ws.Range("I84").Value = 1294.5555
$ws.Range("J84").Value = 1500
$ws.Range("K84").Value = 12945.555
$ws.Range("L84").Value = 15000
$ws.Range("M84").Value = -7641.555
$ws.Range("N84").Value = -25608
